$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.035.77'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '2.374.85'
$ws.Range('E3').Value = '  +6.57%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').Value = "'322.91"
$ws.Range('E5').Value = '  +9.41%  '
$ws.Range('D6').Value = "'103.54"
$ws.Range('E6').Value = '  -8.00%  '
$ws.Range('D7').Value = "'0.643"
$ws.Range('E7').Value = '  +2.50%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = "'0.656"
$ws.Range('E9').Value = '  +9.43%  '
$ws.Range('D10').Value = "'41.22"
$ws.Range('E10').Value = '  -5.25%  '
$ws.Range('D11').Value = "'0.0937"
$ws.Range('E11').Value = '  +2.28%  '
$ws.Range('D12').Value = "'8.53"
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('E13').Value = '  -3.50%  '
$ws.Range('D14').Value = "'16.84"
$ws.Range('E14').Value = '  +12.35%  '
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('D16').Value = '2.733.86'
$ws.Range('E16').Value = '  +6.65%  '
$ws.Range('D17').Value = '2.371.65'
$ws.Range('E17').Value = '  +6.22%  '
$ws.Range('D18').Value = '43.009.39'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').Value = "'7.92"
$ws.Range('E19').Value = '  +9.55%  '
$ws.Range('E20').Value = '  +2.76%  '
$ws.Range('D21').Value = "'76.72"
$ws.Range('E21').Value = '  +4.26%  '
$ws.Range('D22').Value = "'276.90"
$ws.Range('E22').Value = '  +15.57%  '
$ws.Range('E23').Value = '  -0.55%  '
$ws.Range('D25').Value = "'9.49"
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = "'11.67"
$ws.Range('E27').Value = '  +1.93%  '
$ws.Range('D28').Value = "'23.19"
$ws.Range('E28').Value = '  +7.23%  '
$ws.Range('D29').Value = "'175.40"
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = "'37.63"
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('D32').Value = "'3.15"
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').Value = "'0.0916"
$ws.Range('E33').Value = '  +3.95%  '
$ws.Range('D34').Value = "'5.84"
$ws.Range('E34').Value = '  +2.52%  '
$ws.Range('E35').Value = '  +5.08%  '
$ws.Range('E36').Value = '  -1.59%  '
$ws.Range('D37').Value = "'4.15"
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('E38').Value = '  -2.70%  '
$ws.Range('E39').Value = '  +1.72%  '
$ws.Range('D40').Value = "'2.82"
$ws.Range('E40').Value = '  +17.82%  '
$ws.Range('D41').Value = "'1.59"
$ws.Range('E41').Value = '  +21.78%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = "'0.230"
$ws.Range('E42').Value = '  +1.01%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = "'123.48"
$ws.Range('E43').Value = '  +20.79%  '
$ws.Range('D44').Value = "'69.39"
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').Value = "'94.57"
$ws.Range('E45').Value = '  +64.28%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = "'1.00"
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').Value = "'12.37"
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('E48').Value = '  +11.31%  '
$ws.Range('D49').Value = "'5.55"
$ws.Range('E49').Value = '  +2.47%  '
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('D51').Value = '1.607.25'
$ws.Range('E51').Value = '  +12.41%  '
